$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.427.86"
Set-TextValue $ws.Range("E2") "  -1.31%  "
Set-TextValue $ws.Range("D3") "1.847.39"
Set-TextValue $ws.Range("E3") "  -0.99%  "
Set-TextValue $ws.Range("E4") "  -0.02%  "
Set-TextValue $ws.Range("D5") "264.87"
Set-TextValue $ws.Range("E5") "  -1.62%  "
Set-TextValue $ws.Range("E6") "  -0.03%  "
Set-TextValue $ws.Range("D7") "0.5209"
Set-TextValue $ws.Range("E7") "  -1.55%  "
Set-TextValue $ws.Range("D8") "0.3275"
Set-TextValue $ws.Range("E8") "  -2.41%  "
Set-TextValue $ws.Range("E9") "  +0.40%  "
Set-TextValue $ws.Range("D10") "18.85"
Set-TextValue $ws.Range("E10") "  -4.20%  "
Set-TextValue $ws.Range("E11") "  -0.92%  "
Set-TextValue $ws.Range("D12") "0.07780"
Set-TextValue $ws.Range("E12") "  +0.36%  "
Set-TextValue $ws.Range("D13") "1.814.29"
Set-TextValue $ws.Range("E13") "  -2.76%  "
Set-TextValue $ws.Range("D14") "88.19"
Set-TextValue $ws.Range("E14") "  -1.69%  "
Set-TextValue $ws.Range("E15") "  -1.76%  "
Set-TextValue $ws.Range("E16") "  -0.09%  "
Set-TextValue $ws.Range("E17") "  -3.00%  "
Set-TextValue $ws.Range("E18") "  -0.15%  "
Set-TextValue $ws.Range("D19") "0.9997"
Set-TextValue $ws.Range("E19") "  -0.11%  "
Set-TextValue $ws.Range("D20") "26.455.46"
Set-TextValue $ws.Range("E20") "  -1.25%  "
Set-TextValue $ws.Range("D21") "2.074.70"
Set-TextValue $ws.Range("E21") "  -0.74%  "
Set-TextValue $ws.Range("D22") "4.647"
Set-TextValue $ws.Range("E22") "  -0.17%  "
Set-TextValue $ws.Range("D23") "9.553"
Set-TextValue $ws.Range("E23") "  -3.30%  "
Set-TextValue $ws.Range("D24") "6.003"
Set-TextValue $ws.Range("E24") "  -0.83%  "
Set-TextValue $ws.Range("D25") "144.50"
Set-TextValue $ws.Range("E25") "  -0.73%  "
Set-TextValue $ws.Range("D26") "2.200"
Set-TextValue $ws.Range("E26") "  -8.14%  "
Set-TextValue $ws.Range("D27") "1.667"
Set-TextValue $ws.Range("E27") "  +0.82%  "
Set-TextValue $ws.Range("D28") "17.00"
Set-TextValue $ws.Range("E28") "  -1.06%  "
Set-TextValue $ws.Range("D29") "112.01"
Set-TextValue $ws.Range("E29") "  -0.70%  "
Set-TextValue $ws.Range("E30") "  -2.78%  "
Set-TextValue $ws.Range("D31") "4.141"
Set-TextValue $ws.Range("E31") "  -3.01%  "
Set-TextValue $ws.Range("D32") "0.08756"
Set-TextValue $ws.Range("E32") "  -0.89%  "
Set-TextValue $ws.Range("D33") "0.04840"
Set-TextValue $ws.Range("E33") "  -2.01%  "
Set-TextValue $ws.Range("D34") "1.136"
Set-TextValue $ws.Range("E34") "  -1.77%  "
Set-TextValue $ws.Range("D35") "0.7198"
Set-TextValue $ws.Range("D36") "2.851"
Set-TextValue $ws.Range("E36") "  -1.02%  "
Set-TextValue $ws.Range("D37") "3.098"
Set-TextValue $ws.Range("E37") "  -2.56%  "
Set-TextValue $ws.Range("D38") "0.01781"
Set-TextValue $ws.Range("E38") "  -2.85%  "
Set-TextValue $ws.Range("D39") "2.211"
Set-TextValue $ws.Range("E39") "  -3.75%  "
Set-TextValue $ws.Range("D40") "0.4869"
Set-TextValue $ws.Range("E40") "  -3.44%  "
Set-TextValue $ws.Range("D41") "0.9132"
Set-TextValue $ws.Range("E41") "  +2.01%  "
Set-TextValue $ws.Range("D42") "111.40"
Set-TextValue $ws.Range("D43") "6.068"
Set-TextValue $ws.Range("E43") "  -0.82%  "
Set-TextValue $ws.Range("D44") "1.000"
Set-TextValue $ws.Range("E44") "  -0.04%  "
Set-TextValue $ws.Range("D45") "7.716"
Set-TextValue $ws.Range("E45") "  -2.77%  "
Set-TextValue $ws.Range("D46") "0.05939"
Set-TextValue $ws.Range("E46") "  +0.08%  "
Set-TextValue $ws.Range("D47") "0.4167"
Set-TextValue $ws.Range("E47") "  -4.63%  "
Set-TextValue $ws.Range("D48") "9.077"
Set-TextValue $ws.Range("E48") "  -1.91%  "
Set-TextValue $ws.Range("E49") "  -5.93%  "
Set-TextValue $ws.Range("D50") "35.02"
Set-TextValue $ws.Range("E50") "  -2.39%  "
Set-TextValue $ws.Range("D51") "0.8927"
Set-TextValue $ws.Range("E51") "  +2.12%  "
